# "Have TableView working with moveable cards."
#
# Fill in the previously-blank Time Log row 78 (Sheet1) with a new
# time entry: 2014-10-07, 23:10 -> 00:25 (+1 day), 20 min interruption,
# activity "Coding". Row 78's E column already carries the shared
# formula (si="0") inherited from E4, so simply populating A/B/C/D/F
# makes Excel compute E78 and ripple the change through the "Total
# Time:" sum (E83) and the Sheet2 SUMIF/percentage table + pie chart
# that are driven off Sheet1.
#
# Note: D78 (the "Interruption" minutes) is written before B78/C78 so
# the delta formula in E78 - which reads A/B/C/D - recalculates against
# the final set of inputs instead of an intermediate state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A78").Value = 41919                   # Date: 10/7/2014
$ws.Range("D78").Value = 20                       # Interruption: 20 mins
$ws.Range("B78").Value = 0.96527777777777779      # Start Time: 11:10 PM
$ws.Range("C78").Value = 1.0173611111111112       # Stop Time: 12:25 AM (next day)
$ws.Range("F78").Value = "Coding"                 # Activity

# Matches the author's final cursor position after entering the row.
$ws.Range("D79").Select()
